# The workbook has a header row (row 1) containing the title
# "sample_q_short" merged across A1:B1. Remove this header line so the
# "#"/"Question" row becomes the new first row and all data shifts up.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(1).Delete()
